# Applies the coin price / 1h-volume refresh captured in the commit
# "Updated cryptos list on Tue Sep 26 15:52:52 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (the default/unstyled look used by every data cell in
# columns D/E) so that forcing text-storage via NumberFormat does not
# leave a stray style applied to the edited cells.
$refStyle = $ws.Range("B2").Style

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $refStyle
}

Set-TextValue "D2" "26.131.59"
Set-TextValue "E2" "  -0.49%  "
Set-TextValue "D3" "1.583.44"
Set-TextValue "E3" "  -0.18%  "
Set-TextValue "D5" "211.17"
Set-TextValue "E5" "  +0.98%  "
Set-TextValue "E6" "  +0.09%  "
Set-TextValue "E7" "  +0.12%  "
Set-TextValue "E8" "  -0.23%  "
Set-TextValue "E9" "  -1.01%  "
Set-TextValue "D10" "19.17"
Set-TextValue "E10" "  -2.25%  "
Set-TextValue "D11" "0.0846"
Set-TextValue "E11" "  +0.17%  "
Set-TextValue "D12" "1.806.97"
Set-TextValue "E12" "  -0.05%  "
Set-TextValue "D13" "1.611.85"
Set-TextValue "E13" "  +1.59%  "
Set-TextValue "E14" "  -1.69%  "
Set-TextValue "D15" "0.516"
Set-TextValue "E15" "  -0.21%  "
Set-TextValue "D16" "63.97"
Set-TextValue "E16" "  -1.18%  "
Set-TextValue "D17" "26.171.47"
Set-TextValue "E17" "  -0.34%  "
Set-TextValue "E18" "  -0.60%  "
Set-TextValue "E19" "  -0.99%  "
Set-TextValue "D20" "213.07"
Set-TextValue "E20" "  +0.17%  "
Set-TextValue "E21" "  +0.06%  "
Set-TextValue "E22" "  -0.75%  "
Set-TextValue "E23" "  -0.50%  "
Set-TextValue "D24" "8.92"
Set-TextValue "E24" "  +0.53%  "
Set-TextValue "D25" "143.78"
Set-TextValue "E25" "  -0.62%  "
Set-TextValue "E26" "  +0.11%  "
Set-TextValue "D27" "6.97"
Set-TextValue "E27" "  -0.98%  "
Set-TextValue "D29" "15.11"
Set-TextValue "E29" "  -1.35%  "
Set-TextValue "E30" "  -2.24%  "
Set-TextValue "E31" "  +0.31%  "
Set-TextValue "D32" "3.18"
Set-TextValue "E32" "  -1.56%  "
Set-TextValue "D33" "1.338.44"
Set-TextValue "E33" "  +3.86%  "
Set-TextValue "E34" "  -2.13%  "
Set-TextValue "E35" "  +0.04%  "
Set-TextValue "E36" "  -1.50%  "
Set-TextValue "E37" "  -4.21%  "
Set-TextValue "E38" "  -0.24%  "
Set-TextValue "E39" "  +0.23%  "
Set-TextValue "E40" "  +2.53%  "
Set-TextValue "E41" "  +0.05%  "
Set-TextValue "D42" "0.941"
Set-TextValue "E42" "  -17.30%  "
Set-TextValue "D43" "0.766"
Set-TextValue "E43" "  +0.48%  "
Set-TextValue "E44" "  -0.28%  "
Set-TextValue "D45" "1.719.62"
Set-TextValue "E45" "  +0.00%  "
Set-TextValue "D46" "60.87"
Set-TextValue "E46" "  -2.87%  "
Set-TextValue "D47" "85.93"
Set-TextValue "E47" "  -3.18%  "
Set-TextValue "E48" "  -2.18%  "
Set-TextValue "E49" "  -1.76%  "
Set-TextValue "E50" "  -1.11%  "
Set-TextValue "D51" "0.999"
Set-TextValue "E51" "  +0.00%  "
